$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prix Spot")

# Copy formatting from the last existing header cell (AE1) onto the new header cell AF1
$ws.Range("AE1").Copy()
$ws.Range("AF1").PasteSpecial(-4122)
$ws.Range("AF1").Value = "15-jul"

$ws.Range("AF2").Value = 76.77
$ws.Range("AF3").Value = 57.2
$ws.Range("AF4").Value = 63.01
$ws.Range("AF5").Value = 59.29
$ws.Range("AF6").Value = 54.1
$ws.Range("AF7").Value = 63.78
$ws.Range("AF8").Value = 69.52
$ws.Range("AF9").Value = 84.47
$ws.Range("AF10").Value = 94.16
$ws.Range("AF11").Value = 72.97
$ws.Range("AF12").Value = 30.5
$ws.Range("AF13").Value = 39.87
$ws.Range("AF14").Value = 46.2
$ws.Range("AF15").Value = 34.82
$ws.Range("AF16").Value = 25.8
$ws.Range("AF17").Value = 21.39
$ws.Range("AF18").Value = 20.82
$ws.Range("AF19").Value = 38.4
$ws.Range("AF20").Value = 44.77
$ws.Range("AF21").Value = 65.89
$ws.Range("AF22").Value = 74.45999999999999
$ws.Range("AF23").Value = 79.56999999999999
$ws.Range("AF24").Value = 100.79
$ws.Range("AF25").Value = 78.19
